# Revert "Example Data 2 Amelia Tristan"
# Restores the "Data" sheet's D/E columns to the original Eye color / Waist
# variables (undoing the Inseam / Hair Color edit), and restores the
# "Codebook" sheet's corresponding documentation rows.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsCode = $wb.Worksheets.Item("Codebook")

# --- Data sheet: header row -----------------------------------------------
$wsData.Range("D1").Value = "Eye color"
$wsData.Range("E1").Value = "Waist"

# --- Data sheet: replace the Inseam (D) / Hair Color (E) values with the
# original Eye color (D) / Waist (E) values -------------------------------
$wsData.Range("D2").Value = "Blue"
$wsData.Range("E2").Value = 36

$wsData.Range("D3").Value = "Green"
$wsData.Range("E3").Value = 25

$wsData.Range("D4").Value = "Brown"
$wsData.Range("E4").Value = 28

$wsData.Range("D5").Value = "Brown"
$wsData.Range("E5").Value = 35

$wsData.Range("D6").Value = "Black"
$wsData.Range("E6").Value = 36

$wsData.Range("D7").Value = "White"
$wsData.Range("E7").Value = 30

$wsData.Range("D8").Value = "Indigo"
$wsData.Range("E8").Value = 27

$wsData.Range("D9").Value = "Blue"
$wsData.Range("E9").Value = 37

$wsData.Range("D10").Value = "Brown"
$wsData.Range("E10").Value = 40

$wsData.Range("D11").Value = "Green"
$wsData.Range("E11").Value = 32

$wsData.Range("D12").Value = "Black"
$wsData.Range("E12").Value = 48

$wsData.Range("D13").Value = "White"
$wsData.Range("E13").Value = 42

$wsData.Range("D14").Value = "Indigo"
$wsData.Range("E14").Value = 36

$wsData.Range("D15").Value = "Brown"
$wsData.Range("E15").Value = 38

# --- Codebook sheet: add back the Eye Color / Waist documentation rows ----
$wsCode.Range("A5").Value = "Eye Color"
$wsCode.Range("B5").Value = "Natural eye color"
$wsCode.Range("C5").Value = "Black brown, blue"

$wsCode.Range("A6").Value = "Waist"
$wsCode.Range("B6").Value = "Waist incentimeters"
$wsCode.Range("C6").Value = "numeric value >0 or NA"

# --- View / selection state -------------------------------------------
# Data sheet: scroll down and select F20:F21 (no longer the active tab)
$wsData.Activate()
$wsData.Range("F20:F21").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1

# Codebook sheet becomes the active/selected tab, with D6 selected
$wsCode.Activate()
$wsCode.Range("D6").Select()

$wb.Save()
